# "new data to sheet 4"
# Adds an "NFA$millions" column (H) to Sheet4, copied from the same data
# already present in Sheet2's K column, and updates the sheet/window
# selection state so Sheet4 becomes the active tab.

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item(2)
$ws4 = $wb.Worksheets.Item(4)

# --- Sheet4: new column H ("NFA$millions") -------------------------------

# Header cell, styled like the other header cells in row 1 (bold).
$ws4.Range("H1").Value = 'NFA$millions'
$ws4.Range("H1").Font.Bold = $true

# Data values H2:H62 mirror Sheet2's K2:K62 ("NFA$millions") column.
$src = $ws2.Range("K2:K62").Value2
$ws4.Range("H2:H62").Value = $src

# Size the new column similarly to its neighbours.
$ws4.Columns("H:H").ColumnWidth = 12

# Make sure the printed page keeps its (portrait) orientation now that the
# sheet has more content.
$ws4.PageSetup.Orientation = 1

# --- View / selection state -----------------------------------------------

# Sheet2 keeps its own selection but moves off K1:K62 after reviewing it,
# and stops being the focused tab.
$ws2.Activate()
$ws2.Range("K1:K62").Select()
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 39

# Sheet4 becomes the active / focused sheet, with K10 selected.
$ws4.Activate()
$ws4.Range("K10").Select()
